$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style (bold, centered, bordered) from an existing header cell (H1) to I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$values = @(8,8,7,7,6,7,9,8,8,7,8,8,8,7,9,9,8,7,8,7,7,8,8,9,8,9,7,5,8,8,11,9,7,6,8,6,8,7,7,9,7,7,9,8,7,9,8,10,7,8,8,6,8,8,8,9,8,5,8,8,8,8,7,8,8,8,8,8,7,8,9,7,6,4,6)

for ($r = 0; $r -lt $values.Length; $r++) {
    $rowNum = $r + 2
    $v = $values[$r]
    $ws.Cells.Item($rowNum, 9).Value = $v
    $ws.Cells.Item($rowNum, 10).Value = $v
}
